$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (only columns that change: E,G,H,I,J,K,M,N,O,P,Q,R,S,T)
$data = @{
    2  = @{ E=3; G=12.59642866666667;  H=37.789286;   I=0.07012550280485508; J=0.07012550280485508; K=3; M=106.5625623333333; N=319.687687; O=0.4373345410925676; P=0.4373345410925676; Q=1342.307714969054;  R=12080.76943472148;  S=0.03066830458804686; T=0.03066830458804686 }
    3  = @{ E=3; G=12.59642866666667;  H=37.789286;   I=0.07012550280485508; J=0.07012550280485508; K=3; M=102.9000496666667; N=308.700149; O=0.4223035277493257; P=0.4223035277493257; Q=1296.173135422624;  R=11665.55821880362;  S=0.02961424721968554; T=0.02961424721968554 }
    4  = @{ E=3; G=12.59642866666667;  H=37.789286;   I=0.07012550280485508; J=0.07012550280485508; K=3; M=34.20111066666666;  N=102.603332; O=0.1403619311581067; P=0.1403619311581067; Q=430.8118508334391;   R=3877.306657500952;  S=0.009842950997122686; T=0.009842950997122686 }
    5  = @{ E=3; G=107.4733173333333; H=322.419952;  I=0.5983140631002458;  J=0.5983140631002458;  K=3; M=106.5625623333333; N=319.687687; O=0.4373345410925676; P=0.4373345410925676; Q=11452.63207750345;  R=103073.688697531;   S=0.2616634062151755;  T=0.2616634062151755 }
    6  = @{ E=3; G=107.4733173333333; H=322.419952;  I=0.5983140631002458;  J=0.5983140631002458;  K=3; M=102.9000496666667; N=308.700149; O=0.4223035277493257; P=0.4223035277493257; Q=11059.00969144143;  R=99531.08722297284;  S=0.2526701395492665;  T=0.2526701395492665 }
    7  = @{ E=3; G=107.4733173333333; H=322.419952;  I=0.5983140631002458;  J=0.5983140631002458;  K=3; M=34.20111066666666;  N=102.603332; O=0.1403619311581067; P=0.1403619311581067; Q=3675.706819831118;  R=33081.36137848006;  S=0.08398051733580381; T=0.08398051733580381 }
    8  = @{ E=3; G=59.55718233333334;  H=178.671547;  I=0.3315604340948992;  J=0.3315604340948992;  K=3; M=106.5625623333333; N=319.687687; O=0.4373345410925676; P=0.4373345410925676; Q=6346.565954793532;  R=57119.09359314178;  S=0.1450028302893452;  T=0.1450028302893452 }
    9  = @{ E=3; G=59.55718233333334;  H=178.671547;  I=0.3315604340948992;  J=0.3315604340948992;  K=3; M=102.9000496666667; N=308.700149; O=0.4223035277493257; P=0.4223035277493257; Q=6128.437020106723;  R=55155.93318096051;  S=0.1400191409803737;  T=0.1400191409803737 }
    10 = @{ E=3; G=59.55718233333334;  H=178.671547;  I=0.3315604340948992;  J=0.3315604340948992;  K=3; M=34.20111066666666;  N=102.603332; O=0.1403619311581067; P=0.1403619311581067; Q=2036.921783977178;  R=18332.2960557946;   S=0.0465384628251802;  T=0.0465384628251802 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Range("E$r").Value = $row.E
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("M$r").Value = $row.M
    $ws.Range("N$r").Value = $row.N
    $ws.Range("O$r").Value = $row.O
    $ws.Range("P$r").Value = $row.P
    $ws.Range("Q$r").Value = $row.Q
    $ws.Range("R$r").Value = $row.R
    $ws.Range("S$r").Value = $row.S
    $ws.Range("T$r").Value = $row.T
}

$wb.Save()
